$d = $word.ActiveDocument

function Set-RunFormat($rng, $bold) {
    $rng.Font.Name = "Palatino Linotype"
    $rng.Font.Size = 10
    $rng.Bold = $bold
}

# --- Change 1: "June 19, 2022" -> "June 20, 2022" (arraignment date) ---
$d.Content.Find.Execute("June 19, 2022", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "June 20, 2022", 2) | Out-Null

# --- Change 2: append "and report to jail on September 30, 2022, at 7:00 p.m."
#     right after the existing "September 27, 2022" (fines/costs due date) ---
$r = $d.Content
$r.Find.Execute("by September 27, 2022.") | Out-Null
$r.MoveEnd(1, -1) | Out-Null   # wdCharacter; shrink end by 1 to exclude the trailing "."
$r.Collapse(0) | Out-Null      # wdCollapseEnd -> right after "September 27, 2022"

# 2a: ", " (not bold)
$start = $r.Start
$r.InsertAfter(", ")
$seg = $d.Range($start, $r.End)
$seg.Bold = 0

# 2b: "and report to jail on " (not bold)
$r.Collapse(0) | Out-Null
$start = $r.Start
$r.InsertAfter("and report to jail on ")
$seg = $d.Range($start, $r.End)
$seg.Bold = 0

# 2c: "September 30, 2022" (bold)
$r.Collapse(0) | Out-Null
$start = $r.Start
$r.InsertAfter("September 30, 2022")
$seg = $d.Range($start, $r.End)
$seg.Bold = 1

# 2d: "," (bold)
$r.Collapse(0) | Out-Null
$start = $r.Start
$r.InsertAfter(",")
$seg = $d.Range($start, $r.End)
$seg.Bold = 1

# 2e: " at 7:00 " (bold)
$r.Collapse(0) | Out-Null
$start = $r.Start
$r.InsertAfter(" at 7:00 ")
$seg = $d.Range($start, $r.End)
$seg.Bold = 1

# 2f: "p.m" (bold)
$r.Collapse(0) | Out-Null
$start = $r.Start
$r.InsertAfter("p.m")
$seg = $d.Range($start, $r.End)
$seg.Bold = 1

# --- Change 3: insert a new "Restitution." paragraph block right before
#     the existing "Fines and Costs." paragraph text ---
$r2 = $d.Content
$r2.Find.Execute("Fines and Costs.  ") | Out-Null
$p = $r2.Paragraphs(1)
$r2 = $p.Range
$r2.Collapse(1) | Out-Null   # wdCollapseStart -> very start of that paragraph

# 3a: "Restitution." (bold)
$start = $r2.Start
$r2.InsertAfter("Restitution.")
$seg = $d.Range($start, $r2.End)
Set-RunFormat $seg 1

# 3b: " The Defendant must pay restitution in the amount of " (not bold)
$r2.Collapse(0) | Out-Null
$start = $r2.Start
$r2.InsertAfter(" The Defendant must pay restitution in the amount of ")
$seg = $d.Range($start, $r2.End)
Set-RunFormat $seg 0

# 3c: "$5,000 to Justin Kudela no later than September 27, 2022 in order to
#     successfully complete the diversion program. " followed by two manual
#     line breaks (not bold)
$r2.Collapse(0) | Out-Null
$start = $r2.Start
$lbreak = [char]11
$r2.InsertAfter('$5,000 to Justin Kudela no later than September 27, 2022 in order to successfully complete the diversion program. ' + $lbreak + $lbreak)
$seg = $d.Range($start, $r2.End)
Set-RunFormat $seg 0

Write-Output "edits applied"
